# Notes for owl animal companion.
# Add a new "Owl_AC" worksheet (a copy of the Kestrekel_AC reference sheet,
# the last sheet in the workbook) with updated "Special" feat progression
# entries for the Owl animal companion.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Kestrekel_AC")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Leave the reference sheet's selection where the author last left it.
$src.Activate()
$src.Range("H30").Select()

# Copy the Kestrekel_AC sheet to the end of the workbook and rename it.
$src.Copy($null, $lastSheet)
$owl = $wb.Worksheets.Item($wb.Worksheets.Count)
$owl.Name = "Owl_AC"

# Update the "Special" column (E) entries for the Owl animal companion's
# feat progression table.
$owl.Range("E6").Value = "Alertness*"
$owl.Range("E7").Value = "Imp Nat Att"
$owl.Range("E9").Value = "WF(Creature)"
$owl.Range("E10").Value = "Dodge"
$owl.Range("E12").Value = "Mobility"

$owl.Range("E16").Value = "Spring Attack"
$owl.Range("E18").Value = "Improved Critical"
$owl.Range("E19").Value = "Blind-fight"

# Make the new sheet the active tab with the expected selection, matching
# how it was left after editing in Excel.
$owl.Activate()
$owl.Range("E12").Select()
